$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.383.28"
$ws.Range("E2").Value = "  -3.11%  "

$ws.Range("D3").Value = "2.985.40"
$ws.Range("E3").Value = "  -2.90%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'547.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").Value = "'130.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.63%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "2.980.16"
$ws.Range("E8").Value = "  -2.92%  "

$ws.Range("D10").Value = "'5.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.45%  "

$ws.Range("D11").Value = "'0.143"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.37%  "

$ws.Range("D12").Value = "'0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.89%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'33.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.97%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000217"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.03%  "

$ws.Range("D15").Value = "3.468.50"
$ws.Range("E15").Value = "  -2.94%  "

$ws.Range("D16").Value = "61.562.56"
$ws.Range("E16").Value = "  -2.76%  "

$ws.Range("D17").Value = "'0.110"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.71%  "

$ws.Range("D18").Value = "2.997.91"
$ws.Range("E18").Value = "  -2.50%  "

$ws.Range("D19").Value = "'6.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("D20").Value = "'477.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("D21").Value = "'13.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.71%  "

$ws.Range("E22").Value = "  -5.19%  "

$ws.Range("D23").Value = "'6.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.62%  "

$ws.Range("D24").Value = "'80.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.35%  "

$ws.Range("D25").Value = "'11.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.74%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").Value = "'2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").Value = "'7.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.47%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("D30").Value = "'1.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("D31").Value = "'25.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.90%  "

$ws.Range("E32").Value = "  -3.38%  "

$ws.Range("D33").Value = "'2.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("D34").Value = "'5.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("D35").Value = "'54.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.98%  "

$ws.Range("E36").Value = "  -2.64%  "

$ws.Range("D37").Value = "'447.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.07%  "

$ws.Range("D38").Value = "3.120.56"
$ws.Range("E38").Value = "  -4.30%  "

$ws.Range("E39").Value = "  -0.81%  "

$ws.Range("E40").Value = "  -5.89%  "

$ws.Range("D41").Value = "'0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.19%  "

$ws.Range("E42").Value = "  -0.92%  "

$ws.Range("E44").Value = "  -10.55%  "

$ws.Range("D45").Value = "'25.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("E46").Value = "  -4.91%  "

$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("D48").Value = "'1.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.14%  "

$ws.Range("D50").Value = "'113.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.64%  "

$ws.Range("D51").Value = "0.0₃0479"
$ws.Range("E51").Value = "  -9.24%  "
